$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 20 - "取指定成员参数" (getMemberParam) now documents the "id"
#    key inside the json param sample.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = '"{"id":[],"param":["id",]}"'

# ---------------------------------------------------------------------
# 2) New "部门" (branch/department) rows 36-38, under the existing
#    "获取部门成员" row (row 35). Fill in interface url, params,
#    param description and meaning for each new api.
# ---------------------------------------------------------------------
$ws.Range("B36").Value = "http://xxx/abranch!getPosition"
$ws.Range("C36").Value = "无"
$ws.Range("D36").Value = "无"
$ws.Range("E36").Value = "获取职位"

$ws.Range("B37").Value = "http://xxx/abranch!getBranchTreeAndMember"
$ws.Range("C37").Value = "无"
$ws.Range("D37").Value = "无"
$ws.Range("E37").Value = "获取部门及成员"

$ws.Range("B38").Value = "http://xxx/getBranchTree"
$ws.Range("C38").Value = "无"
$ws.Range("D38").Value = "无"
$ws.Range("E38").Value = "获取部门树"

# Hyperlink the new interface-url cells, same as every other "接口" cell.
$ws.Hyperlinks.Add($ws.Range("B36"), "http://xxx/abranch!getPosition")
$ws.Hyperlinks.Add($ws.Range("B37"), "http://xxx/abranch!getBranchTreeAndMember")
$ws.Hyperlinks.Add($ws.Range("B38"), "http://xxx/getBranchTree")

# Re-apply the B column "interface" style (blue/underline hyperlink look)
# that Hyperlinks.Add doesn't reuse from the existing style table.
$ws.Range("B35").Copy()
$ws.Range("B36:B38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Merge the module column A35:A38 so "部门" spans every branch row,
#    same pattern as the other module groups above it.
# ---------------------------------------------------------------------
$ws.Range("A35:A38").Merge()

# ---------------------------------------------------------------------
# 4) Restore the view: scrolled down a bit further, with B44 selected.
# ---------------------------------------------------------------------
$ws.Range("B44").Select()

Write-Output "edit complete"
